$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: convert B3 from text "15-01-22" to a real date value with date
# number format (built-in id 14, "mm-dd-yy") and left alignment, then fill
# in the remark/location/engineer columns.
$ws.Range("B3").HorizontalAlignment = -4131
$ws.Range("B3").NumberFormat = "mm-dd-yy"
$ws.Range("B3").Value = 44689
$ws.Range("C3").Value = "size not updated"
$ws.Range("E3").Value = "CCM "
$ws.Range("F3").Value = "kunal"

# Row 4: same breakdown record, reuse B3's format via copy/paste so the
# style table doesn't grow a duplicate entry.
$ws.Range("B3").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("B4").Value = 44689
$ws.Range("C4").Value = "size not updated"
$ws.Range("E4").Value = "CCM "
$ws.Range("F4").Value = "kunal"

# Update the current selection to match the saved workbook state.
$ws.Range("G4").Select()
